# Daily attendance processing - 2025-10-25 08:46:23
#
# Re-orders the "Recorded By" e-mail lists for several sessions, refreshes
# attendance counters / coverage percentages, fills in the previously
# "Not Recorded" Year2/B3/HISTOLOGY session (row 97), tweaks the Status
# column width, and re-records the Year2/B4/HISTOLOGY attendance numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: write a literal text value into a cell without Excel's "smart"
# re-interpretation of numeric/percentage-looking strings (e.g. "11.7%"
# would otherwise be silently converted into the number 0.117 with a
# percentage number format, which also touches the cell's style). We
# build the literal string as a formula result on a scratch cell, copy
# it, and paste-special just the VALUE into the destination - this keeps
# the destination's existing style/number-format completely untouched.
# ---------------------------------------------------------------------
function Set-LiteralText {
    param($range, [string]$text)

    $scratch = $ws.Range("ZZ1")
    $escaped = $text.Replace('"', '""')
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
    $scratch.ClearContents()
}

# ---------------------------------------------------------------------
# "Recorded By" e-mail list re-ordering (same addresses, new order)
# ---------------------------------------------------------------------
$anatomyList = "nahla.nagiub@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("G2").Value = $anatomyList
$ws.Range("G17").Value = $anatomyList
$ws.Range("G92").Value = $anatomyList
$ws.Range("G107").Value = $anatomyList

$histologyList = "Arwa.elnagar@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"
$ws.Range("G7").Value = $histologyList
$ws.Range("G22").Value = $histologyList
$ws.Range("G112").Value = $histologyList

$anatomyListA3A4 = "nahla.nagiub@med.asu.edu.eg, gehanadel@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G32").Value = $anatomyListA3A4
$ws.Range("G47").Value = $anatomyListA3A4

$anatomyListB1B2 = "nahla.nagiub@med.asu.edu.eg, gehanadel@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G62").Value = $anatomyListB1B2
$ws.Range("G77").Value = $anatomyListB1B2

# ---------------------------------------------------------------------
# Attendance counters / student fractions
# ---------------------------------------------------------------------
$ws.Range("H22").Value = "117/217"
$ws.Range("H112").Value = "97/226"

# ---------------------------------------------------------------------
# Class Statistics block (columns K:S)
# ---------------------------------------------------------------------
$ws.Range("L6").Value = 14          # Recorded Sessions
$ws.Range("L7").Value = 0           # Missing Sessions
Set-LiteralText $ws.Range("L9") "11.7%"   # Coverage %
Set-LiteralText $ws.Range("L10") "54.6%"  # Average Attendance %

Set-LiteralText $ws.Range("S16") "56.7%"  # Year2 / A2 attendance %

$ws.Range("O21").Value = 2          # Year2 / B3 Absent
$ws.Range("P21").Value = 0          # Year2 / B3 Late
Set-LiteralText $ws.Range("R21") "13.3%"  # Year2 / B3 absence %
Set-LiteralText $ws.Range("S21") "36.6%"  # Year2 / B3 attendance %

Set-LiteralText $ws.Range("S22") "43.8%"  # Year2 / B4 attendance %

# ---------------------------------------------------------------------
# Row 97 (Year2 / B3 / HISTOLOGY) - session got recorded: copy the
# formatting used by already-recorded rows (style index 2, as seen on
# row 2) over the previously "Not Recorded" style (index 9), then fill
# in the recorder e-mails, the attendance fraction and the new status.
# ---------------------------------------------------------------------
$ws.Range("A2:I2").Copy()
$ws.Range("A97:I97").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("G97").Value = "Arwa.elnagar@med.asu.edu.eg, aya.saeed@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"
$ws.Range("H97").Value = "1/224"
$ws.Range("I97").Value = "Recorded"

# ---------------------------------------------------------------------
# Narrow the "Status" column (column I / 9th column)
# ---------------------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth = 9.16666666666667

$excel.CutCopyMode = 0
